$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.911.33"
$ws.Range("D2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.888.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7735"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.80"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3101"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.62"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07158"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08575"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7636"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.913.43"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.361"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.50%  "

$ws.Range("E15").Value = "  +1.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.146"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.896.54"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.76"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.07"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007810"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.189.69"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9976"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.949"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1640"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.348"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.96"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.74"
$ws.Range("D28").ClearFormats()

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.032"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.440"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.09%  "

$ws.Range("E31").Value = "  -1.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.509"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.104"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("E34").Value = "  -1.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.240"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7464"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("E38").Value = "  +2.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01960"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.782"
$ws.Range("D40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4465"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.107.53"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.95%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.15"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.60%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.080"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8505"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.868"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.609"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.091.28"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.990"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.12%  "
